$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Reset all "Value" (grading) entries to 0 ahead of the announcement
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0

$ws.Range("G8").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("G12").Value = 0

$ws.Range("G15").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("G21").Value = 0

# Update Unit weights for the "Calling Services" criteria rows
$ws.Range("D24").Value = 2
$ws.Range("D25").Value = 2
$ws.Range("D26").Value = 2

$ws.Range("G24").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("G28").Value = 0

# Move the active selection to the Extra subtotal cell
[void]$ws.Range("G29").Select()
